$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: corrections to five existing cells -----------------------
# (Q56, Q62: a stray isPivot-week flag gets cleared; O988 picks up a
#  restated pivot count; R990/R991 "backup" placeholders get their
#  real value of 0 now that the row is finalised.)
$ws.Cells.Item(56, 17).Value  = 0     # Q56  : 1 -> 0
$ws.Cells.Item(62, 17).Value  = 0     # Q62  : 1 -> 0
$ws.Cells.Item(988, 15).Value = 2     # O988 : 0 -> 2
$ws.Cells.Item(990, 18).Value = 0     # R990 : (blank) -> 0
$ws.Cells.Item(991, 18).Value = 0     # R991 : (blank) -> 0

# --- Part 2: append 25 new weekly OHLCV rows (992-1016) ---------------
# Columns: A Datetime(serial) B Open C High D Low E Close F AdjClose
#          G Volume H Year I Month J Day K Hour L Minute M Second
#          N Week O isPivot P two_line_structure Q detect_structure R backup
# F (Adj Close) and R (backup) are not yet populated for these freshly
# appended rows, so they stay blank - same as the source feed.
$newRows = @(
    @(45474, 23.78000068664551, 27.04000091552734, 23.61000061035156, 26.63999938964844, $null, 1511973132, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0, $null),
    @(45481, 26.95000076293945, 27.44000053405762, 24.60000038146973, 25.60000038146973, $null, 2153691920, 2024, 7, 8, 0, 0, 0, 28, 1, 0, 0, $null),
    @(45488, 25.64999961853027, 26.85000038146973, 24.6200008392334, 24.77000045776367, $null, 913604801, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1, $null),
    @(45495, 25.89999961853027, 26.19000053405762, 24.25, 24.96999931335449, $null, 878292252, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0, $null),
    @(45502, 25.25, 27.20000076293945, 24.79999923706055, 25.45000076293945, $null, 1574024328, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0, $null),
    @(45509, 24.5, 24.75, 23.15999984741211, 23.94000053405762, $null, 948897467, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0, $null),
    @(45516, 23.88999938964844, 25.10000038146973, 23.54000091552734, 24.27000045776367, $null, 586574800, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0, $null),
    @(45523, 24.29000091552734, 24.94000053405762, 24.25, 24.38999938964844, $null, 498804833, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 2, $null),
    @(45530, 24.5, 24.53000068664551, 23.54999923706055, 23.61000061035156, $null, 421883396, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0, $null),
    @(45537, 23.61000061035156, 24.02000045776367, 23.10000038146973, 23.13999938964844, $null, 427939155, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0, $null),
    @(45544, 23.04999923706055, 24.40999984741211, 22.51000022888184, 23.29999923706055, $null, 636180013, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0, $null),
    @(45551, 23.30999946594238, 23.96999931335449, 22.79999923706055, 23, $null, 470073033, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0, $null),
    @(45558, 23.19000053405762, 23.19000053405762, 22.64999961853027, 22.73999977111816, $null, 304299480, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0, $null),
    @(45565, 22.72999954223633, 22.73999977111816, 21.52000045776367, 21.95000076293945, $null, 347338972, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0, $null),
    @(45572, 22.20999908447266, 22.26000022888184, 20.77000045776367, 21.39999961853027, $null, 446878765, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0, $null),
    @(45579, 21.39999961853027, 21.48999977111816, 20.56999969482422, 21.09000015258789, $null, 250419774, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0, $null),
    @(45586, 21.10000038146973, 21.52000045776367, 19.21999931335449, 19.3799991607666, $null, 395300821, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0, $null),
    @(45593, 20.35000038146973, 21.29000091552734, 20.10000038146973, 20.70999908447266, $null, 485903164, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0, $null),
    @(45600, 20.70000076293945, 21.13999938964844, 20.06999969482422, 20.15999984741211, $null, 292924603, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0, $null),
    @(45607, 20.15999984741211, 20.32999992370605, 19.02000045776367, 19.30999946594238, $null, 274366308, 2024, 11, 11, 0, 0, 0, 46, 2, 0, 0, $null),
    @(45614, 19.29000091552734, 19.81999969482422, 19.04000091552734, 19.20999908447266, $null, 223370399, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0, $null),
    @(45621, 19.51000022888184, 20.64999961853027, 19.05999946594238, 19.95999908447266, $null, 611320063, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0, $null),
    @(45628, 19.95999908447266, 21.57999992370605, 19.89999961853027, 21.48999977111816, $null, 564295134, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0, $null),
    @(45635, 21.48999977111816, 22.02000045776367, 20.70000076293945, 21.09000015258789, $null, 548915070, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0, $null),
    @(45642, 21.14999961853027, 21.57999992370605, 19.75, 19.82999992370605, $null, 434170151, 2024, 12, 16, 0, 0, 0, 51, 0, 0, 0, $null)
)

$startRow = 992
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value  = $row[1]   # Open
    $ws.Cells.Item($r, 3).Value  = $row[2]   # High
    $ws.Cells.Item($r, 4).Value  = $row[3]   # Low
    $ws.Cells.Item($r, 5).Value  = $row[4]   # Close

    # Adj Close: blank placeholder (quote-prefixed empty text, then
    # restyled to Normal so no stray quote-prefix formatting lingers)
    $adjClose = $ws.Cells.Item($r, 6)
    $adjClose.Value = "'"
    $adjClose.Style = "Normal"

    $ws.Cells.Item($r, 7).Value  = $row[6]   # Volume
    $ws.Cells.Item($r, 8).Value  = $row[7]   # Year
    $ws.Cells.Item($r, 9).Value  = $row[8]   # Month
    $ws.Cells.Item($r, 10).Value = $row[9]   # Day
    $ws.Cells.Item($r, 11).Value = $row[10]  # Hour
    $ws.Cells.Item($r, 12).Value = $row[11]  # Minute
    $ws.Cells.Item($r, 13).Value = $row[12]  # Second
    $ws.Cells.Item($r, 14).Value = $row[13]  # Week
    $ws.Cells.Item($r, 15).Value = $row[14]  # isPivot
    $ws.Cells.Item($r, 16).Value = $row[15]  # two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[16]  # detect_structure

    # backup: same blank-placeholder treatment as Adj Close
    $backup = $ws.Cells.Item($r, 18)
    $backup.Value = "'"
    $backup.Style = "Normal"
}
